# ---------------------------------------------------------------------------
# Add 2022-Q1 data:
#   1. The existing "总计" sheet (aggregate totals) is repurposed to hold the
#      per-fund holdings detail for the new 2022-Q1 quarter (it becomes the
#      "2022-Q1" sheet), mirroring how each prior quarter got its own sheet.
#   2. A brand-new "总计" sheet is appended right after it, containing the
#      same aggregate table as before plus a new leading row for 2022-Q1.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# Template cells to clone formatting from (bold+bordered header / index style).
$tplSheet = $wb.Worksheets.Item("2020-Q4")
$tplHeaderRow = $tplSheet.Range("B1:H1")
$tplIndexCol  = $tplSheet.Range("A2:A13")

# --- Step 1: repurpose "总计" -> "2022-Q1" ---------------------------------
$ws6 = $wb.Worksheets.Item("总计")
$ws6.Cells.ClearContents()
$ws6.Name = "2022-Q1"

# Header row B1:H1 + index column A2:A12 get the shared bold/border style.
$tplHeaderRow.Copy()
$ws6.Range("B1:H1").PasteSpecial(-4122)
$tplIndexCol.Copy()
$ws6.Range("A2:A12").PasteSpecial(-4122)

# Force the numeric-looking text columns (fund code / scale / position pct /
# ratio pct / market value) to stay text, matching the source data which
# stores these as formatted strings rather than numbers.
$ws6.Range("B2:G12").NumberFormat = "@"

    $ws6.Range("B1").Value = "基金代码"
    $ws6.Range("C1").Value = "基金名称"
    $ws6.Range("D1").Value = "基金规模"
    $ws6.Range("E1").Value = "股票总仓位"
    $ws6.Range("F1").Value = "仓位占比"
    $ws6.Range("G1").Value = "持有市值(亿元)"
    $ws6.Range("H1").Value = "仓位排名"
    $ws6.Range("A2").Value = 0
    $ws6.Range("B2").Value = "159828"
    $ws6.Range("C2").Value = "国泰中证医疗ETF"
    $ws6.Range("D2").Value = "13.16"
    $ws6.Range("E2").Value = "99.03"
    $ws6.Range("F2").Value = "3.39"
    $ws6.Range("G2").Value = "0.4461"
    $ws6.Range("H2").Value = 9
    $ws6.Range("A3").Value = 1
    $ws6.Range("B3").Value = "004634"
    $ws6.Range("C3").Value = "新疆前海联合泳涛灵活配置混合A"
    $ws6.Range("D3").Value = "1.33"
    $ws6.Range("E3").Value = "89.65"
    $ws6.Range("F3").Value = "4.93"
    $ws6.Range("G3").Value = "0.0656"
    $ws6.Range("H3").Value = 5
    $ws6.Range("A4").Value = 2
    $ws6.Range("B4").Value = "002310"
    $ws6.Range("C4").Value = "创金合信沪深300指数增强A"
    $ws6.Range("D4").Value = "3.71"
    $ws6.Range("E4").Value = "91.95"
    $ws6.Range("F4").Value = "1.55"
    $ws6.Range("G4").Value = "0.0575"
    $ws6.Range("H4").Value = 10
    $ws6.Range("A5").Value = 3
    $ws6.Range("B5").Value = "004050"
    $ws6.Range("C5").Value = "华夏新锦升灵活配置混合A"
    $ws6.Range("D5").Value = "1.35"
    $ws6.Range("E5").Value = "69.58"
    $ws6.Range("F5").Value = "3.70"
    $ws6.Range("G5").Value = "0.0500"
    $ws6.Range("H5").Value = 9
    $ws6.Range("A6").Value = 4
    $ws6.Range("B6").Value = "006235"
    $ws6.Range("C6").Value = "东方城镇消费主题混合"
    $ws6.Range("D6").Value = "0.50"
    $ws6.Range("E6").Value = "90.32"
    $ws6.Range("F6").Value = "7.23"
    $ws6.Range("G6").Value = "0.0362"
    $ws6.Range("H6").Value = 3
    $ws6.Range("A7").Value = 5
    $ws6.Range("B7").Value = "002315"
    $ws6.Range("C7").Value = "创金合信沪深300指数增强C"
    $ws6.Range("D7").Value = "2.29"
    $ws6.Range("E7").Value = "91.95"
    $ws6.Range("F7").Value = "1.55"
    $ws6.Range("G7").Value = "0.0355"
    $ws6.Range("H7").Value = 10
    $ws6.Range("A8").Value = 6
    $ws6.Range("B8").Value = "011002"
    $ws6.Range("C8").Value = "同泰大健康主题混合A"
    $ws6.Range("D8").Value = "0.44"
    $ws6.Range("E8").Value = "90.75"
    $ws6.Range("F8").Value = "5.55"
    $ws6.Range("G8").Value = "0.0244"
    $ws6.Range("H8").Value = 7
    $ws6.Range("A9").Value = 7
    $ws6.Range("B9").Value = "011003"
    $ws6.Range("C9").Value = "同泰大健康主题混合C"
    $ws6.Range("D9").Value = "0.30"
    $ws6.Range("E9").Value = "90.75"
    $ws6.Range("F9").Value = "5.55"
    $ws6.Range("G9").Value = "0.0166"
    $ws6.Range("H9").Value = 7
    $ws6.Range("A10").Value = 8
    $ws6.Range("B10").Value = "004135"
    $ws6.Range("C10").Value = "申万菱信量化成长混合"
    $ws6.Range("D10").Value = "0.49"
    $ws6.Range("E10").Value = "86.91"
    $ws6.Range("F10").Value = "1.98"
    $ws6.Range("G10").Value = "0.0097"
    $ws6.Range("H10").Value = 5
    $ws6.Range("A11").Value = 9
    $ws6.Range("B11").Value = "004051"
    $ws6.Range("C11").Value = "华夏新锦升灵活配置混合C"
    $ws6.Range("D11").Value = "0.00"
    $ws6.Range("E11").Value = "69.58"
    $ws6.Range("F11").Value = "3.70"
    $ws6.Range("G11").Value = 0
    $ws6.Range("H11").Value = 9
    $ws6.Range("A12").Value = 10
    $ws6.Range("B12").Value = "007041"
    $ws6.Range("C12").Value = "新疆前海联合泳涛灵活配置混合C"
    $ws6.Range("D12").Value = "0.00"
    $ws6.Range("E12").Value = "89.65"
    $ws6.Range("F12").Value = "4.93"
    $ws6.Range("G12").Value = 0
    $ws6.Range("H12").Value = 5

# --- Step 2: add the new "总计" sheet right after "2022-Q1" ----------------
$ws7 = $wb.Worksheets.Add($null, $ws6)
$ws7.Name = "总计"

# Header row B1:D1 + index column A2:A7 get the shared bold/border style.
$tplHeaderRow.Copy()
$ws7.Range("B1:D1").PasteSpecial(-4122)
$tplIndexCol.Copy()
$ws7.Range("A2:A7").PasteSpecial(-4122)

    $ws7.Range("B1").Value = "日期"
    $ws7.Range("C1").Value = "持有数量(只)"
    $ws7.Range("D1").Value = "持有市值(亿元)"
    $ws7.Range("A2").Value = 0
    $ws7.Range("B2").Value = "2022-Q1"
    $ws7.Range("C2").Value = 11
    $ws7.Range("D2").Value = 0.74
    $ws7.Range("A3").Value = 1
    $ws7.Range("B3").Value = "2021-Q4"
    $ws7.Range("C3").Value = 11
    $ws7.Range("D3").Value = 0.8100000000000001
    $ws7.Range("A4").Value = 2
    $ws7.Range("B4").Value = "2021-Q3"
    $ws7.Range("C4").Value = 39
    $ws7.Range("D4").Value = 6.04
    $ws7.Range("A5").Value = 3
    $ws7.Range("B5").Value = "2021-Q2"
    $ws7.Range("C5").Value = 80
    $ws7.Range("D5").Value = 21.94
    $ws7.Range("A6").Value = 4
    $ws7.Range("B6").Value = "2021-Q1"
    $ws7.Range("C6").Value = 21
    $ws7.Range("D6").Value = 4.69
    $ws7.Range("A7").Value = 5
    $ws7.Range("B7").Value = "2020-Q4"
    $ws7.Range("C7").Value = 6
    $ws7.Range("D7").Value = 3.33
